$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 100
$ws.Range("I9").Value = 100
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 100
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = 69
$ws.Range("N9").ClearContents()

$ws.Range("H17").Value = 1877.7778
$ws.Range("J17").Value = 1877.7778
$ws.Range("L17").Value = 5633.3334
$ws.Range("N17").Value = -5969.3334

$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()

$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("M47").ClearContents()

$ws.Range("H58").Value = 1105.7693
$ws.Range("I58").Value = 257.5
$ws.Range("J58").Value = 3933.3333
$ws.Range("K58").Value = 772.5
$ws.Range("L58").Value = 11799.9999
$ws.Range("M58").Value = -622.5
$ws.Range("N58").Value = -12099.9999

$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()

$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()

$ws.Range("H112").Value = 3834.2144
$ws.Range("J112").Value = 3548.25
$ws.Range("L112").Value = 10644.75
$ws.Range("N112").Value = -12860.75

$ws.Range("H129").Value = 1462.3846
$ws.Range("I129").Value = 841.1
$ws.Range("J129").Value = 3533.3333
$ws.Range("K129").Value = 2523.3
$ws.Range("L129").Value = 10599.9999
$ws.Range("M129").Value = 2476.7
$ws.Range("N129").Value = -20599.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1841.625
$ws.Range("I45").Value = 1544.4546
$ws.Range("K45").Value = 1544.4546
$ws.Range("M45").Value = -1167.4546

$ws.Range("H74").Value = 3600.4
$ws.Range("J74").Value = 4500
$ws.Range("L74").Value = 4500
$ws.Range("N74").Value = -6248

$ws.Range("H77").Value = 3600.4
$ws.Range("J77").Value = 4500
$ws.Range("L77").Value = 22500
$ws.Range("N77").Value = -31236

$ws.Range("H102").Value = 27780328
$ws.Range("I102").Value = 27780328
$ws.Range("K102").Value = 27780328
$ws.Range("M102").Value = -27778706

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4449.75
$ws.Range("I20").Value = 3639.8
$ws.Range("K20").Value = 3639.8
$ws.Range("M20").Value = -3392.8

$ws.Range("H22").Value = 200
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()

$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("N46").ClearContents()

$ws.Range("H64").Value = 944.5714
$ws.Range("J64").Value = 992.4
$ws.Range("L64").Value = 992.4
$ws.Range("N64").Value = -1442.4

$ws.Range("H67").Value = 944.5714
$ws.Range("J67").Value = 992.4
$ws.Range("L67").Value = 992.4
$ws.Range("N67").Value = -2552.4

$ws.Range("H105").Value = 6197.6665
$ws.Range("I105").Value = 4943
$ws.Range("K105").Value = 4943
$ws.Range("M105").Value = -3196

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 121322.4
$ws.Range("I22").Value = 171071.28
$ws.Range("J22").Value = 5241.6665
$ws.Range("K22").Value = 171071.28
$ws.Range("L22").Value = 5241.6665
$ws.Range("M22").Value = -170721.28
$ws.Range("N22").Value = -5941.6665

$ws.Range("H31").Value = 9330.556
$ws.Range("I31").Value = 7995
$ws.Range("K31").Value = 7995
$ws.Range("M31").Value = -7700

$ws.Range("H34").Value = 9330.556
$ws.Range("I34").Value = 7995
$ws.Range("K34").Value = 7995
$ws.Range("M34").Value = -7793

$ws.Range("H56").Value = 3093
$ws.Range("I56").Value = 3093
$ws.Range("K56").Value = 3093
$ws.Range("M56").Value = -2248

$ws.Range("H94").Value = 1824.3846
$ws.Range("I94").Value = 1145.2
$ws.Range("K94").Value = 1145.2
$ws.Range("M94").Value = -694.2

$ws.Range("H105").Value = 601.7222
$ws.Range("I105").Value = 622.13336
$ws.Range("K105").Value = 622.13336
$ws.Range("M105").Value = 1124.86664

$ws.Range("H134").Value = 2374.7334
$ws.Range("J134").Value = 3207.8333
$ws.Range("L134").Value = 9623.499899999999
$ws.Range("N134").Value = -14693.4999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 2055.262
$ws.Range("I11").Value = 2217
$ws.Range("K11").Value = 6651
$ws.Range("M11").Value = -6511

$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("M44").ClearContents()
$ws.Range("N44").ClearContents()

$ws.Range("H57").Value = 8665.666999999999
$ws.Range("I57").Value = 8249
$ws.Range("K57").Value = 24747
$ws.Range("M57").Value = -24188

$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").ClearContents()

$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 15062.625
$ws.Range("I70").Value = 12099.2
$ws.Range("K70").Value = 12099.2
$ws.Range("M70").Value = -11829.2

$ws.Range("H73").Value = 15062.625
$ws.Range("I73").Value = 12099.2
$ws.Range("K73").Value = 12099.2
$ws.Range("M73").Value = -11163.2

$ws.Range("H132").Value = 3249.3333
$ws.Range("J132").Value = 3499.6667
$ws.Range("L132").Value = 10499.0001
$ws.Range("N132").Value = -15559.0001

$ws.Range("H136").Value = 138666.33
$ws.Range("J136").Value = 138666.33
$ws.Range("L136").Value = 415998.99
$ws.Range("N136").Value = -421098.99

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4967.3335
$ws.Range("I7").Value = 4967.3335
$ws.Range("K7").Value = 4967.3335
$ws.Range("M7").Value = -4855.3335

$ws.Range("H22").Value = 1276
$ws.Range("I22").Value = 1047.25
$ws.Range("K22").Value = 1047.25
$ws.Range("M22").Value = -752.25

$ws.Range("H27").Value = 1276
$ws.Range("I27").Value = 1047.25
$ws.Range("K27").Value = 1047.25
$ws.Range("M27").Value = -940.25

$ws.Range("H34").Value = 49500
$ws.Range("I34").Value = 49000
$ws.Range("J34").Value = 50000
$ws.Range("K34").Value = 49000
$ws.Range("L34").Value = 50000
$ws.Range("M34").Value = -48828
$ws.Range("N34").Value = -50344

$ws.Range("H46").Value = 2714
$ws.Range("I46").Value = 2879.6
$ws.Range("K46").Value = 2879.6
$ws.Range("M46").Value = -2691.6

$ws.Range("H68").Value = 9686.875
$ws.Range("I68").Value = 3332.8333
$ws.Range("K68").Value = 3332.8333
$ws.Range("M68").Value = -2583.8333

$ws.Range("H71").Value = 9686.875
$ws.Range("I71").Value = 3332.8333
$ws.Range("K71").Value = 16664.1665
$ws.Range("M71").Value = -12920.1665

$ws.Range("H82").Value = 85389
$ws.Range("I82").Value = 2242.6365
$ws.Range("K82").Value = 2242.6365
$ws.Range("M82").Value = -1881.6365

$ws.Range("H85").Value = 85389
$ws.Range("I85").Value = 2242.6365
$ws.Range("K85").Value = 2242.6365
$ws.Range("M85").Value = -994.6365000000001

$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()

$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()

$ws.Range("H126").Value = 4967.3335
$ws.Range("I126").Value = 4967.3335
$ws.Range("K126").Value = 14902.0005
$ws.Range("M126").Value = -12432.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 1817
$ws.Range("I8").Value = 101.5
$ws.Range("J8").Value = 2674.75
$ws.Range("K8").Value = 101.5
$ws.Range("L8").Value = 2674.75
$ws.Range("M8").Value = 38.5
$ws.Range("N8").Value = -2954.75

$ws.Range("H24").Value = 3000
$ws.Range("I24").Value = 3000
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 3000
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = -2770
$ws.Range("N24").ClearContents()

$ws.Range("H74").Value = 54812
$ws.Range("J74").Value = 54812
$ws.Range("L74").Value = 54812
$ws.Range("N74").Value = -56684

$ws.Range("H77").Value = 54812
$ws.Range("J77").Value = 54812
$ws.Range("L77").Value = 164436
$ws.Range("N77").Value = -173796

$ws.Range("H94").Value = 26947
$ws.Range("I94").Value = 26947
$ws.Range("K94").Value = 26947
$ws.Range("M94").Value = -26046
